$d = $word.ActiveDocument

# Locate the two paragraphs that still contain the "clientRegion" template
# placeholder (they are textually identical, so each must be processed
# individually, in document order).
$targets = @()
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*clientRegion*") {
        $targets += $idx
    }
}

$findText = "} обл., `${clientTown},  `${clientStreet}, буд., `${clientHouse}, `${clientFlat}"

foreach ($pIndex in $targets) {
    $rng = $d.Paragraphs($pIndex).Range
    # wdReplaceOne (1) so only the first match at/after the start of this
    # paragraph's range is touched; the other, textually identical
    # paragraph is left alone during this iteration. wdFindStop (0) keeps
    # the search confined to the supplied range.
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, "} ", 1) | Out-Null
}

# The first (table) paragraph keeps a trailing space after the closing
# brace ("} "), while the second paragraph's closing brace has no trailing
# space ("}"). Trim that extra space back off on the last paragraph only.
$lastIndex = $targets[$targets.Count - 1]
$rng2 = $d.Paragraphs($lastIndex).Range
# Exclude the trailing paragraph mark from the search range so the literal
# "} " search matches against the run text rather than "} " + pilcrow.
$rng2.End = $rng2.End - 1
$rng2.Find.Execute("} ", $true, $false, $false, $false, $false, $true, 0, $false, "}", 1) | Out-Null
